$wb = $excel.ActiveWorkbook

# --- "VerifyStartEndDateValidation" sheet: test result now fails ---
$ws = $wb.Worksheets.Item("VerifyStartEndDateValidation")

# Order matters for how new shared strings get appended, so set the
# "expected message" text before the "actual user" text.
$ws.Range("I2").Value = "End Dat should always be greater or equal to the Start Dat!"
$ws.Range("H2").Value = "'Test User"
$ws.Range("L2").Value = "expected [End Dat should always be greater or equal to the Start Dat!] but found []"
$ws.Range("K2").Value = "FAIL"
$ws.Range("A2").Value = "Test Manager"

# Reset view: scroll back to the top-left and select D20.
$ws.Activate()
$ws.Range("D20").Select()

# --- "Test Cases" summary sheet: reflect the new FAIL result ---
$ws1 = $wb.Worksheets.Item("Test Cases")
$ws1.Range("F2").Value = "FAIL"
